$wb = $excel.ActiveWorkbook

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 3
$ws.Range("H3").Value = 17899.8
$ws.Range("I3").Value = 9876
$ws.Range("J3").Value = 49995
$ws.Range("K3").Value = 9876
$ws.Range("L3").Value = 49995
$ws.Range("M3").Value = -9761
$ws.Range("N3").Value = -50225
# Row 15
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
# Row 24
$ws.Range("H24").Value = 46998.5
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 46998.5
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 46998.5
$ws.Range("N24").Value = -47746.5
# Row 32
$ws.Range("H32").Value = 5082.2
$ws.Range("I32").Value = 4155.8374
$ws.Range("J32").Value = 24999
$ws.Range("K32").Value = 4155.8374
$ws.Range("L32").Value = 24999
$ws.Range("M32").Value = -3868.8374
$ws.Range("N32").Value = -25573
# Row 35
$ws.Range("H35").Value = 1300
$ws.Range("I35").Value = 1300
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1300
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -894
# Row 61
$ws.Range("H61").Value = 4463.885
$ws.Range("I61").Value = 3928.25
$ws.Range("J61").Value = 5320.9
$ws.Range("K61").Value = 3928.25
$ws.Range("L61").Value = 5320.9
$ws.Range("M61").Value = -3716.25
$ws.Range("N61").Value = -5744.9
# Row 74
$ws.Range("H74").Value = 30306548
$ws.Range("I74").Value = 41668068
$ws.Range("J74").Value = 9154.666999999999
$ws.Range("K74").Value = 41668068
$ws.Range("L74").Value = 9154.666999999999
$ws.Range("M74").Value = -41667194
$ws.Range("N74").Value = -10902.667
# Row 77
$ws.Range("H77").Value = 30306548
$ws.Range("I77").Value = 41668068
$ws.Range("J77").Value = 9154.666999999999
$ws.Range("K77").Value = 208340340
$ws.Range("L77").Value = 45773.335
$ws.Range("M77").Value = -208335972
$ws.Range("N77").Value = -54509.335
# Row 100
$ws.Range("H100").Value = 46998.5
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 46998.5
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 46998.5
$ws.Range("N100").Value = -49162.5
# Row 132
$ws.Range("H132").Value = 2734.675
$ws.Range("I132").Value = 2143.9644
$ws.Range("J132").Value = 4113
$ws.Range("K132").Value = 6431.8932
$ws.Range("L132").Value = 12339
$ws.Range("M132").Value = -3901.8932
$ws.Range("N132").Value = -17399
# Row 136
$ws.Range("H136").Value = 4463.885
$ws.Range("I136").Value = 3928.25
$ws.Range("J136").Value = 5320.9
$ws.Range("K136").Value = 11784.75
$ws.Range("L136").Value = 15962.7
$ws.Range("M136").Value = -9234.75
$ws.Range("N136").Value = -21062.7

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 10
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
# Row 11
$ws.Range("H11").Value = 549.75
$ws.Range("I11").Value = 550
$ws.Range("J11").Value = 549.5
$ws.Range("K11").Value = 550
$ws.Range("L11").Value = 549.5
$ws.Range("M11").Value = -410
$ws.Range("N11").Value = -829.5
# Row 20
$ws.Range("H20").Value = 4659.7036
$ws.Range("I20").Value = 4742.5557
$ws.Range("J20").Value = 4494
$ws.Range("K20").Value = 4742.5557
$ws.Range("L20").Value = 4494
$ws.Range("M20").Value = -4495.5557
$ws.Range("N20").Value = -4988
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
# Row 86
$ws.Range("H86").Value = 4986.5625
$ws.Range("I86").Value = 4349.154
$ws.Range("J86").Value = 7748.6665
$ws.Range("K86").Value = 4349.154
$ws.Range("L86").Value = 7748.6665
$ws.Range("M86").Value = -3226.154
$ws.Range("N86").Value = -9994.666499999999
# Row 89
$ws.Range("H89").Value = 4986.5625
$ws.Range("I89").Value = 4349.154
$ws.Range("J89").Value = 7748.6665
$ws.Range("K89").Value = 21745.77
$ws.Range("L89").Value = 38743.3325
$ws.Range("M89").Value = -16129.77
$ws.Range("N89").Value = -49975.3325
# Row 134
$ws.Range("H134").Value = 2826.4783
$ws.Range("I134").Value = 1814.7858
$ws.Range("J134").Value = 4400.222
$ws.Range("K134").Value = 5444.357400000001
$ws.Range("L134").Value = 13200.666
$ws.Range("M134").Value = -2909.357400000001
$ws.Range("N134").Value = -18270.666

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 3007.4102
$ws.Range("I132").Value = 2156.8064
$ws.Range("J132").Value = 6303.5
$ws.Range("K132").Value = 6470.4192
$ws.Range("L132").Value = 18910.5
$ws.Range("M132").Value = -3940.4192
$ws.Range("N132").Value = -23970.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 114194.22
$ws.Range("I11").Value = 146321.14
$ws.Range("J11").Value = 1750
$ws.Range("K11").Value = 438963.42
$ws.Range("L11").Value = 5250
$ws.Range("M11").Value = -438823.42
$ws.Range("N11").Value = -5530
# Row 37
$ws.Range("H37").Value = 311358
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 311358
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 934074
$ws.Range("N37").Value = -934298
# Row 68
$ws.Range("H68").Value = 1659.2727
$ws.Range("I68").Value = 948.2
$ws.Range("J68").Value = 2251.8333
$ws.Range("K68").Value = 2844.6
$ws.Range("L68").Value = 6755.499899999999
$ws.Range("M68").Value = -2033.6
$ws.Range("N68").Value = -8377.499899999999
# Row 71
$ws.Range("H71").Value = 1659.2727
$ws.Range("I71").Value = 948.2
$ws.Range("J71").Value = 2251.8333
$ws.Range("K71").Value = 8533.800000000001
$ws.Range("L71").Value = 20266.4997
$ws.Range("M71").Value = -4477.800000000001
$ws.Range("N71").Value = -28378.4997
# Row 122
$ws.Range("H122").Value = 2876.8235
$ws.Range("I122").Value = 873.3333
$ws.Range("J122").Value = 3306.1428
$ws.Range("K122").Value = 7859.9997
$ws.Range("L122").Value = 29755.2852
$ws.Range("M122").Value = -5409.9997
$ws.Range("N122").Value = -34655.2852

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 123
$ws.Range("H123").Value = 35011.145
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 35011.145
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 35011.145
$ws.Range("N123").Value = -39911.145
# Row 133
$ws.Range("H133").Value = 72000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 72000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 72000
$ws.Range("N133").Value = -82120
# Row 135
$ws.Range("H135").Value = 68394.086
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 68394.086
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 68394.086
$ws.Range("N135").Value = -78534.086
# Row 140
$ws.Range("H140").Value = 70072
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 70072
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 70072
$ws.Range("N140").Value = -80432

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 9
$ws.Range("H9").Value = 198.33333
$ws.Range("I9").Value = 247.5
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 247.5
$ws.Range("L9").Value = 100
$ws.Range("M9").Value = -23.5
$ws.Range("N9").Value = -548
# Row 132
$ws.Range("H132").Value = 3509.7659
$ws.Range("I132").Value = 3149.2856
$ws.Range("J132").Value = 4041
$ws.Range("K132").Value = 9447.856800000001
$ws.Range("L132").Value = 12123
$ws.Range("M132").Value = -6917.856800000001
$ws.Range("N132").Value = -17183

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 330.83334
$ws.Range("I113").Value = 319.05554
$ws.Range("J113").Value = 366.16666
$ws.Range("K113").Value = 957.16662
$ws.Range("L113").Value = 1098.49998
$ws.Range("M113").Value = 1212.83338
$ws.Range("N113").Value = -5438.499980000001
